$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the existing row 77, shifting rows 77:92 down to 80:95
$ws.Rows("77:79").Insert()

# Row 77 - new record (Castle Brite / Especial)
$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 44889
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100103
$ws.Range("H77").Value = "Frutos de hueso (carozo)"
$ws.Range("I77").Value = 100103003
$ws.Range("J77").Value = "Damasco"
$ws.Range("K77").Value = "Castle Brite"
$ws.Range("L77").Value = "Especial"
$ws.Range("M77").Value = 380
$ws.Range("N77").Value = 30000
$ws.Range("O77").Value = 30000
$ws.Range("P77").Value = 30000
$ws.Range("Q77").Value = "$/caja 15 kilos granel"
$ws.Range("R77").Value = "Provincia de Los Andes"
$ws.Range("S77").Value = 2000
$ws.Range("T77").Value = 15

# Row 78 - new record (Castle Brite / Primera, 15 kilos granel)
$ws.Range("A78").Value = 9
$ws.Range("B78").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C78").Value = "Metropolitana"
$ws.Range("D78").Value = 44889
$ws.Range("E78").Value = 13
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100103
$ws.Range("H78").Value = "Frutos de hueso (carozo)"
$ws.Range("I78").Value = 100103003
$ws.Range("J78").Value = "Damasco"
$ws.Range("K78").Value = "Castle Brite"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 410
$ws.Range("N78").Value = 24000
$ws.Range("O78").Value = 24000
$ws.Range("P78").Value = 24000
$ws.Range("Q78").Value = "$/caja 15 kilos granel"
$ws.Range("R78").Value = "Provincia de Los Andes"
$ws.Range("S78").Value = 1600
$ws.Range("T78").Value = 15

# Row 79 - new record (Castle Brite / Primera, 18 kilos granel)
$ws.Range("A79").Value = 9
$ws.Range("B79").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C79").Value = "Metropolitana"
$ws.Range("D79").Value = 44889
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100103
$ws.Range("H79").Value = "Frutos de hueso (carozo)"
$ws.Range("I79").Value = 100103003
$ws.Range("J79").Value = "Damasco"
$ws.Range("K79").Value = "Castle Brite"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 350
$ws.Range("N79").Value = 28800
$ws.Range("O79").Value = 28800
$ws.Range("P79").Value = 28800
$ws.Range("Q79").Value = "$/caja 18 kilos granel"
$ws.Range("R79").Value = "Provincia de Los Andes"
$ws.Range("S79").Value = 1600
$ws.Range("T79").Value = 18
